# Auto commit at 2025-10-12  7:37:20.44
# Append two new daily rows (2025-10-11) for 四方坪站 and 高岭站 to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New date serial 45941 == 2025-10-11 (zero out time-of-day so the stored
# serial is a clean integer, matching the source workbook's date cells)
$newDate = Get-Date -Year 2025 -Month 10 -Day 11 -Hour 0 -Minute 0 -Second 0

# ---- Row 82: 四方坪站充电量(kw) ----
$row82 = 82
$ws.Cells.Item($row82, 1).Value = $newDate
$ws.Cells.Item($row82, 1).NumberFormat = "yyyy\-mm\-dd"

$ws.Cells.Item($row82, 2).Value = "四方坪站充电量(kw)"

$row82Values = @(
    650.91600000000005,
    1805.6980000000001,
    614.50699999999995,
    289.21699999999998,
    230.96999999999997,
    774.67300000000012,
    687.78899999999999,
    255.61199999999997,
    154.262,
    219.44499999999999,
    228.96799999999999,
    290.70000000000005,
    1086.5380000000002,
    1170.1410000000001,
    681.84699999999987,
    284.46600000000001,
    320.90800000000002,
    271.601,
    85.941999999999993,
    168.36,
    73.893000000000001,
    95.699999999999989,
    35.200000000000003,
    30.768000000000001
)

for ($i = 0; $i -lt $row82Values.Length; $i++) {
    $col = 3 + $i  # Column C = 3
    $cell = $ws.Cells.Item($row82, $col)
    $cell.Value = $row82Values[$i]
    $cell.NumberFormat = "0.00_);[Red]\(0.00\)"
}

# ---- Row 83: 高岭站充电量(kw) ----
$row83 = 83
$ws.Cells.Item($row83, 1).Value = $newDate
$ws.Cells.Item($row83, 1).NumberFormat = "yyyy\-mm\-dd"

$ws.Cells.Item($row83, 2).Value = "高岭站充电量(kw)"

$row83Values = @(
    470.464,
    402.50300000000004,
    283.34100000000001,
    112.515,
    126.26299999999999,
    93.281000000000006,
    163.25900000000001,
    258.25299999999999,
    604.54799999999989,
    210.62699999999998,
    106.637,
    222.51999999999998,
    488.00599999999997,
    599.27499999999986,
    324.86799999999999,
    296.42699999999996,
    152.18,
    95.593999999999994,
    128,
    0,
    0,
    41.198999999999998,
    0,
    11.922000000000001
)

for ($i = 0; $i -lt $row83Values.Length; $i++) {
    $col = 3 + $i  # Column C = 3
    $cell = $ws.Cells.Item($row83, $col)
    $cell.Value = $row83Values[$i]
    $cell.NumberFormat = "0.00_);[Red]\(0.00\)"
}

# ---- Update view state to match the recorded scroll/selection ----
$ws.Application.ActiveWindow.TopLeftCell = $ws.Range("A56")
$ws.Range("D85").Select()
